# Fruta / hortaliza, semanal
# Insert 3 new weekly price rows for "Cereza" (Brooks, Lapins, Rainier)
# ahead of the existing historical rows, which shift down from 97-105 to 100-108.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows at row 97; existing rows 97-105 shift down to 100-108.
$ws.Rows.Item(97).Insert()
$ws.Rows.Item(97).Insert()
$ws.Rows.Item(97).Insert()

# New row 97: Brooks / Primera
$ws.Cells.Item(97,1).Value  = 4
$ws.Cells.Item(97,2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(97,3).Value  = "Los Lagos"
$ws.Cells.Item(97,4).Value  = 44918
$ws.Cells.Item(97,5).Value  = 10
$ws.Cells.Item(97,6).Value  = "Fruta"
$ws.Cells.Item(97,7).Value  = 100103
$ws.Cells.Item(97,8).Value  = "Frutos de hueso (carozo)"
$ws.Cells.Item(97,9).Value  = 100103001
$ws.Cells.Item(97,10).Value = "Cereza"
$ws.Cells.Item(97,11).Value = "Brooks"
$ws.Cells.Item(97,12).Value = "Primera"
$ws.Cells.Item(97,13).Value = 800
$ws.Cells.Item(97,14).Value = 5500
$ws.Cells.Item(97,15).Value = 6000
$ws.Cells.Item(97,16).Value = 5750
$ws.Cells.Item(97,17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(97,18).Value = "Provincia de Curicó"
$ws.Cells.Item(97,19).Value = 575
$ws.Cells.Item(97,20).Value = 10

# New row 98: Lapins / Primera
$ws.Cells.Item(98,1).Value  = 4
$ws.Cells.Item(98,2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(98,3).Value  = "Los Lagos"
$ws.Cells.Item(98,4).Value  = 44918
$ws.Cells.Item(98,5).Value  = 10
$ws.Cells.Item(98,6).Value  = "Fruta"
$ws.Cells.Item(98,7).Value  = 100103
$ws.Cells.Item(98,8).Value  = "Frutos de hueso (carozo)"
$ws.Cells.Item(98,9).Value  = 100103001
$ws.Cells.Item(98,10).Value = "Cereza"
$ws.Cells.Item(98,11).Value = "Lapins"
$ws.Cells.Item(98,12).Value = "Primera"
$ws.Cells.Item(98,13).Value = 1000
$ws.Cells.Item(98,14).Value = 5000
$ws.Cells.Item(98,15).Value = 6000
$ws.Cells.Item(98,16).Value = 5500
$ws.Cells.Item(98,17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(98,18).Value = "Provincia de Curicó"
$ws.Cells.Item(98,19).Value = 550
$ws.Cells.Item(98,20).Value = 10

# New row 99: Rainier / Primera
$ws.Cells.Item(99,1).Value  = 4
$ws.Cells.Item(99,2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(99,3).Value  = "Los Lagos"
$ws.Cells.Item(99,4).Value  = 44918
$ws.Cells.Item(99,5).Value  = 10
$ws.Cells.Item(99,6).Value  = "Fruta"
$ws.Cells.Item(99,7).Value  = 100103
$ws.Cells.Item(99,8).Value  = "Frutos de hueso (carozo)"
$ws.Cells.Item(99,9).Value  = 100103001
$ws.Cells.Item(99,10).Value = "Cereza"
$ws.Cells.Item(99,11).Value = "Rainier"
$ws.Cells.Item(99,12).Value = "Primera"
$ws.Cells.Item(99,13).Value = 1000
$ws.Cells.Item(99,14).Value = 8500
$ws.Cells.Item(99,15).Value = 9000
$ws.Cells.Item(99,16).Value = 8750
$ws.Cells.Item(99,17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(99,18).Value = "Provincia de Curicó"
$ws.Cells.Item(99,19).Value = 875
$ws.Cells.Item(99,20).Value = 10
